# Update "Correspond Handback DateTime" values on the zh-cn and de-de sheets
# to reflect a newly-generated handback report (commit: "Generate Report for Handback").

$wb = $excel.ActiveWorkbook

# zh-cn sheet: Correspond Handoff Datetime (E2) and Correspond Handback DateTime (H2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 11:11:09"
$wsZhCn.Range("H2").Value = "2016-03-22 11:11:37"

# de-de sheet: Correspond Handoff Datetime (E2) and Correspond Handback DateTime (H2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 11:11:13"
$wsDeDe.Range("H2").Value = "2016-03-22 11:11:44"
